$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = "Fri Oct 06 11:13:35 EDT 2023"
$ws.Range("B3").Value = "Fri Oct 06 11:13:46 EDT 2023"
$ws.Range("B4").Value = "Fri Oct 06 11:13:57 EDT 2023"
$ws.Range("B5").Value = "Fri Oct 06 11:14:08 EDT 2023"
$ws.Range("B6").Value = "Fri Oct 06 11:14:18 EDT 2023"
$ws.Range("B7").Value = "Fri Oct 06 11:14:30 EDT 2023"
$ws.Range("B8").Value = "Fri Oct 06 11:14:40 EDT 2023"
$ws.Range("B9").Value = "Fri Oct 06 11:14:51 EDT 2023"
$ws.Range("B10").Value = "Fri Oct 06 11:15:01 EDT 2023"
$ws.Range("B11").Value = "Fri Oct 06 11:15:12 EDT 2023"
$ws.Range("B12").Value = "Fri Oct 06 11:15:24 EDT 2023"
$ws.Range("B13").Value = "Fri Oct 06 11:15:34 EDT 2023"
